$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (India - ISL, Goa vs Punjab) updated odds
$ws.Range("G3").Value = 2.05
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 2.2
$ws.Range("L3").Value = 3.75
$ws.Range("N3").Value = 12
$ws.Range("U3").Value = 1.67
$ws.Range("V3").Value = 2.1
$ws.Range("W3").Value = 8.5
$ws.Range("X3").Value = 11
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 17
$ws.Range("AB3").Value = 26
$ws.Range("AD3").Value = 6.5
$ws.Range("AH3").Value = 11
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 26
$ws.Range("AM3").Value = 29
$ws.Range("AN3").Value = 4.33
$ws.Range("AO3").Value = 11
$ws.Range("AP3").Value = 21
$ws.Range("AQ3").Value = 41
$ws.Range("AX3").Value = 17
$ws.Range("AY3").Value = 23
$ws.Range("AZ3").Value = 51
$ws.Range("BA3").Value = 67

# Row 4 (Indonesia - Liga 1, Persija Jakarta vs Madura United) updated odds
$ws.Range("G4").Value = 1.57
$ws.Range("H4").Value = 3.55
$ws.Range("I4").Value = 5.7
$ws.Range("J4").Value = 2.15
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 5.6
$ws.Range("N4").Value = 7.7
$ws.Range("O4").Value = 1.26
$ws.Range("P4").Value = 3.2
$ws.Range("Q4").Value = 1.82
$ws.Range("R4").Value = 1.9
$ws.Range("U4").Value = 1.75
$ws.Range("V4").Value = 1.85
$ws.Range("X4").Value = 7.5
$ws.Range("Y4").Value = 7.8
$ws.Range("Z4").Value = 12
$ws.Range("AB4").Value = 24
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 70
$ws.Range("AG4").Value = 500
$ws.Range("AH4").Value = 15.5
$ws.Range("AI4").Value = 37
$ws.Range("AK4").Value = 120
$ws.Range("AL4").Value = 60
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 7.8
$ws.Range("AP4").Value = 16.5
$ws.Range("AQ4").Value = 25
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.55
$ws.Range("AU4").Value = 7.3
$ws.Range("AV4").Value = 65
$ws.Range("AW4").Value = 7.1
$ws.Range("BA4").Value = 250
